$d = $word.ActiveDocument
$apos = [char]0x2019

function Replace-In-Paragraph($index, $old, $new) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 2)
}

# Paragraph 1: "English" (hyperlink) + " / Portuguese / French / Thai / Vietnamese / Spanish"
Replace-In-Paragraph 1 "English" "영어"
Replace-In-Paragraph 1 " / Portuguese / French / Thai / Vietnamese / Spanish" " / 포르투갈어 / 프랑스어 / 태국어 / 베트남어 / 스페인어"

# Paragraph 3: standalone "English"
Replace-In-Paragraph 3 "English" "영어"

# Paragraph 5: "Brief" (bold run) + ":"
Replace-In-Paragraph 5 "Brief" "브리핑"

# Paragraph 6: brief description
$p6old = "An email sent to partners in the target country who RSVPed yes but didn" + $apos + "t submit their documents by the deadline. We will be revoking their invites. It will be sent via customer.io"
$p6new = "기한 내에 서류를 제출하지 않은 파트너에게 발송되는 이메일입니다. 초대장을 철회할 예정입니다. It will be sent via customer.io"
Replace-In-Paragraph 6 $p6old $p6new

# Paragraph 8: "Target audience" (bold run) + ":"
Replace-In-Paragraph 8 "Target audience" "대상 청중"

# Paragraph 9: target audience description
$p9old = "Invited partners who didn" + $apos + "t submit their documents on time"
$p9new = "제때 서류를 제출하지 않은 초대된 파트너들"
Replace-In-Paragraph 9 $p9old $p9new

# Paragraph 12: "Subject line" (bold) + ": Your " + [EVENT NAME] + " registration"
Replace-In-Paragraph 12 "Subject line" "제목"
Replace-In-Paragraph 12 ": Your " ": 귀하의 "
Replace-In-Paragraph 12 " registration" " 등록"

# Paragraph 14: headline
$p14old = "We didn" + $apos + "t receive your documents on time"
$p14new = "귀하의 문서를 제시간에 받지 못했습니다"
Replace-In-Paragraph 14 $p14old $p14new

# Paragraph 17: deadline paragraph
$p17old1 = "We didn" + $apos + "t receive your documents by the deadline ("
Replace-In-Paragraph 17 $p17old1 "마감일("
$p17old2 = ")" + ". Unfortunately, we" + $apos + "re unable to proceed with your registration for the "
$p17new2 = ")까지 귀하의 문서를 받지 못했습니다. Unfortunately, we" + $apos + "re unable to proceed with your registration for the "
Replace-In-Paragraph 17 $p17old2 $p17new2

# Paragraph 18: "We wish you..." + commented "conference/seminar/affiliate trip"
Replace-In-Paragraph 18 "We wish you the best and hope to see you at our next " "다음 "
Replace-In-Paragraph 18 "conference/seminar/affiliate trip" "컨퍼런스/세미나/제휴 여행에서 뵙기를 기대합니다."

# Paragraph 20: country manager contact info (NOT paragraph 19, which keeps
# its English "live chat or WhatsApp" text unchanged)
Replace-In-Paragraph 20 "If you have any questions, please contact your country manager, " "궁금하신 사항은, "
Replace-In-Paragraph 20 ", at " "에게 "
Replace-In-Paragraph 20 " or " " 또는 "
Replace-In-Paragraph 20 " (WhatsApp). " " (WhatsApp)을 통해 연락해 주시기 바랍니다. "

# --- Comments: both say "choose either one" -> "하나를 선택하세요" ---
foreach ($c in $d.Comments) {
    $c.Range.Text = "하나를 선택하세요"
}
